$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Table rows (1-indexed, row 1 = header):
#   Row 6 = "SYM"               -> becomes "ARD"
#   Row 7 = "ARD"               -> becomes "Simplified ARD 4" (also gains taller row height)
#   Row 8 = "ER"                -> becomes "SYM"
#   Row 9 = "Simplified ARD 4"  -> becomes "ER" (also gains shorter row height)

# Row 6: SYM -> ARD
$t.Cell(6,1).Range.Text = "ARD"
$t.Cell(6,2).Range.Text = "20"
$t.Cell(6,3).Range.Text = "-2,674.14"
$t.Cell(6,4).Range.Text = "5,388.29"

# Row 7: ARD -> Simplified ARD 4 (row height 676 -> 700 twips = 33.8 -> 35 pt)
$t.Rows.Item(7).Height = 35
$t.Cell(7,1).Range.Text = "Simplified ARD 4"
$t.Cell(7,2).Range.Text = "10"
$t.Cell(7,3).Range.Text = "-2,710.26"
$t.Cell(7,4).Range.Text = "5,440.51"

# Row 8: ER -> SYM
$t.Cell(8,1).Range.Text = "SYM"
$t.Cell(8,2).Range.Text = "10"
$t.Cell(8,3).Range.Text = "-2,780.02"
$t.Cell(8,4).Range.Text = "5,580.03"

# Row 9: Simplified ARD 4 -> ER (row height 700 -> 676 twips = 35 -> 33.8 pt)
$t.Rows.Item(9).Height = 33.8
$t.Cell(9,1).Range.Text = "ER"
$t.Cell(9,2).Range.Text = "1"
$t.Cell(9,3).Range.Text = "-3,256.93"
$t.Cell(9,4).Range.Text = "6,515.85"
